$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1,0.9997,0.9994,0.9991,0.9988,0.9984999999999999,0.9982,0.9979,0.9976,0.9973,0.997,0.9967,0.9964,0.9961,0.9958,0.9955000000000001,0.9952,0.9949,0.9946,0.9943,0.994,0.9937,0.9933999999999999,0.9931,0.9928,0.9924999999999999,0.9922,0.9919,0.9916,0.9913,0.991,0.9907,0.9903999999999999,0.9901,0.9898,0.9895,0.9892,0.9889,0.9886,0.9883,0.988,0.9877,0.9873999999999999,0.9871,0.9868,0.9864999999999999,0.9862,0.9859,0.9856,0.9853,0.985,0.9847,0.9843999999999999,0.9841,0.9838,0.9835,0.9832,0.9829,0.9826,0.9823,0.982,0.9817,0.9813999999999999,0.9811,0.9808,0.9804999999999999,0.9802,0.9799,0.9796,0.9792999999999999,0.979,0.9787,0.9783999999999999,0.9781,0.9778,0.9775,0.9772,0.9769,0.9766,0.9762999999999999,0.976,0.9757,0.9753999999999999,0.9751,0.9748,0.9744999999999999,0.9742,0.9739,0.9736,0.9732999999999999,0.973,0.9727,0.9723999999999999,0.9721,0.9718,0.9715,0.9712,0.9709,0.9706,0.9702999999999999)
for ($i = 0; $i -lt 100; $i++) {
    $ws.Cells.Item(3, $i + 2).Value = $values[$i]
}

$values = @(1,0.9968,0.9936,0.9903999999999999,0.9872,0.984,0.9808,0.9776,0.9743999999999999,0.9712,0.968,0.9648,0.9616,0.9583999999999999,0.9551999999999999,0.952,0.9488,0.9456,0.9423999999999999,0.9391999999999999,0.9359999999999999,0.9328,0.9296,0.9263999999999999,0.9231999999999999)
for ($i = 0; $i -lt 25; $i++) {
    $ws.Cells.Item(4, $i + 2).Value = $values[$i]
}

$values = @(1,0.9984999999999999,0.997,0.9955000000000001,0.994,0.9924999999999999,0.991,0.9894999999999999,0.988,0.9864999999999999,0.985,0.9835,0.982,0.9804999999999999,0.979,0.9775,0.976,0.9744999999999999,0.973,0.9715,0.97,0.9684999999999999,0.967,0.9655,0.964,0.9624999999999999,0.961,0.9595,0.958,0.9564999999999999,0.955,0.9535,0.952,0.9504999999999999,0.949,0.9475,0.946,0.9444999999999999,0.9429999999999999,0.9415,0.9399999999999999,0.9385,0.9369999999999999,0.9355,0.9339999999999999,0.9324999999999999,0.9309999999999999,0.9295,0.9279999999999999,0.9264999999999999)
for ($i = 0; $i -lt 50; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $values[$i]
}

$values = @(1,0.9985000000000001,0.997,0.9955000000000001,0.994,0.9925,0.991,0.9895,0.988,0.9865,0.985,0.9835,0.982,0.9805,0.979,0.9775,0.976,0.9745,0.973,0.9715,0.97,0.9685,0.967,0.9655,0.964,0.9625,0.961,0.9595,0.958,0.9565,0.955,0.9535,0.952,0.9505,0.949,0.9475,0.946,0.9445,0.9429999999999999,0.9415,0.9399999999999999,0.9385,0.9370000000000001,0.9355,0.9339999999999999,0.9325,0.9309999999999999,0.9295,0.9279999999999999,0.9265,0.925,0.9235,0.9219999999999999,0.9205,0.919,0.9175,0.9159999999999999,0.9145,0.913,0.9115,0.91,0.9085,0.907,0.9055,0.904,0.9025,0.901,0.8995,0.898,0.8965,0.895,0.8935,0.892,0.8905,0.889,0.8875,0.886,0.8845,0.883,0.8815,0.88,0.8784999999999999,0.877,0.8754999999999999,0.874,0.8725000000000001,0.871,0.8694999999999999,0.868,0.8664999999999999,0.865,0.8634999999999999,0.862,0.8604999999999999,0.859,0.8574999999999999,0.856,0.8545,0.853,0.8514999999999999)
for ($i = 0; $i -lt 100; $i++) {
    $ws.Cells.Item(7, $i + 2).Value = $values[$i]
}

$values = @(1,0.9928,0.9856,0.9784,0.9712,0.964,0.9568,0.9496,0.9424,0.9352,0.9279999999999999,0.9208,0.9136,0.9064,0.8992,0.892,0.8847999999999999,0.8775999999999999,0.8704,0.8632,0.856,0.8488,0.8415999999999999,0.8343999999999999,0.8271999999999999)
for ($i = 0; $i -lt 25; $i++) {
    $ws.Cells.Item(8, $i + 2).Value = $values[$i]
}

$values = @(1,0.9976,0.9952,0.9928,0.9903999999999999,0.988,0.9856,0.9832,0.9808,0.9783999999999999,0.976,0.9736,0.9712,0.9688,0.9663999999999999,0.964,0.9616,0.9591999999999999,0.9568,0.9543999999999999,0.952,0.9496,0.9471999999999999,0.9448,0.9423999999999999,0.9399999999999999,0.9376,0.9351999999999999,0.9328,0.9303999999999999,0.9279999999999999,0.9256,0.9231999999999999,0.9208,0.9183999999999999,0.9159999999999999,0.9136,0.9111999999999999,0.9087999999999999,0.9063999999999999,0.9039999999999999,0.9016,0.8991999999999999,0.8967999999999999,0.8943999999999999,0.8919999999999999,0.8895999999999999,0.8871999999999999,0.8847999999999999,0.8823999999999999)
for ($i = 0; $i -lt 50; $i++) {
    $ws.Cells.Item(11, $i + 2).Value = $values[$i]
}

$values = @(1,0.9983,0.9966,0.9949,0.9932,0.9915,0.9898,0.9881,0.9864000000000001,0.9847,0.983,0.9813000000000001,0.9796,0.9779,0.9762,0.9745,0.9728,0.9711,0.9694,0.9677,0.966,0.9643,0.9626,0.9609,0.9592000000000001,0.9575,0.9558,0.9541000000000001,0.9524,0.9507,0.9490000000000001,0.9473,0.9456,0.9439000000000001,0.9422,0.9405,0.9388000000000001,0.9371,0.9354,0.9337,0.9320000000000001,0.9303,0.9286000000000001,0.9269000000000001,0.9252,0.9235,0.9218000000000001,0.9201,0.9184,0.9167000000000001)
for ($i = 0; $i -lt 50; $i++) {
    $ws.Cells.Item(13, $i + 2).Value = $values[$i]
}

$values = @(1,0.9972,0.9944,0.9916,0.9888,0.986,0.9832,0.9804,0.9776,0.9748,0.972,0.9692,0.9664,0.9636,0.9608,0.958,0.9552,0.9524,0.9496,0.9468,0.944,0.9412,0.9384,0.9356,0.9328,0.9299999999999999,0.9272,0.9244,0.9216,0.9188,0.916,0.9132,0.9104,0.9076,0.9047999999999999,0.902,0.8992,0.8964,0.8935999999999999,0.8908,0.888,0.8852,0.8824,0.8795999999999999,0.8768,0.874,0.8712,0.8684000000000001,0.8655999999999999,0.8628)
for ($i = 0; $i -lt 50; $i++) {
    $ws.Cells.Item(15, $i + 2).Value = $values[$i]
}
